$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$src = $ws.Range("E4:E10")
$vals = $src.Value()
$dst = $ws.Range("E3:E9")
$dst.Value = $vals
